# Trf-Tfr2.xlsx — refresh with new TPM-based NATMI output.
# The "ECs" cluster is no longer a valid *target* cluster in the refreshed run,
# so every row whose Target cluster (col D) was "ECs" is dropped, shrinking the
# 4x4 sending/target matrix (16 data rows) down to 4x3 (12 data rows), and all
# ligand/receptor expression + specificity figures are recomputed against the
# smaller cluster set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 14:17 (Sending cluster "Resolving-Mac", all four targets)
# up front so the remaining rows 2:13 can simply be overwritten in place with
# the recomputed values below (row counts/order for clusters ECs/FAPs/MuSCs stay
# put; only their D/E..T figures change, and their D="ECs" rows are repointed).
$ws.Rows("14:17").Delete()

# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Trf"
$ws.Range("C2").Value = "Tfr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.600481666666667
$ws.Range("H2").Value = 4.801445
$ws.Range("I2").Value = 0.01283248898116485
$ws.Range("J2").Value = 0.01283248898116485
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.818892666666667
$ws.Range("N2").Value = 8.456678
$ws.Range("O2").Value = 0.9725758139465345
$ws.Range("P2").Value = 0.9725758139465346
$ws.Range("Q2").Value = 4.511586033301112
$ws.Range("R2").Value = 40.60427429971
$ws.Range("S2").Value = 0.01248056841581634
$ws.Range("T2").Value = 0.01248056841581634

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Trf"
$ws.Range("C3").Value = "Tfr2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.600481666666667
$ws.Range("H3").Value = 4.801445
$ws.Range("I3").Value = 0.01283248898116485
$ws.Range("J3").Value = 0.01283248898116485
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04134966666666667
$ws.Range("N3").Value = 0.124049
$ws.Range("O3").Value = 0.01426648349910611
$ws.Range("P3").Value = 0.01426648349910611
$ws.Range("Q3").Value = 0.06617938342277778
$ws.Range("R3").Value = 0.5956144508050001
$ws.Range("S3").Value = 0.0001830744923022493
$ws.Range("T3").Value = 0.0001830744923022493

# Row 4: ECs -> Resolving-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Trf"
$ws.Range("C4").Value = "Tfr2"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.600481666666667
$ws.Range("H4").Value = 4.801445
$ws.Range("I4").Value = 0.01283248898116485
$ws.Range("J4").Value = 0.01283248898116485
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.038136
$ws.Range("N4").Value = 0.114408
$ws.Range("O4").Value = 0.01315770255435942
$ws.Range("P4").Value = 0.01315770255435942
$ws.Range("Q4").Value = 0.06103596884
$ws.Range("R4").Value = 0.54932371956
$ws.Range("S4").Value = 0.0001688460730462618
$ws.Range("T4").Value = 0.0001688460730462618

# Row 5: FAPs -> FAPs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Trf"
$ws.Range("C5").Value = "Tfr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.234235
$ws.Range("H5").Value = 27.702705
$ws.Range("I5").Value = 0.07403909795092109
$ws.Range("J5").Value = 0.07403909795092109
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.818892666666667
$ws.Range("N5").Value = 8.456678
$ws.Range("O5").Value = 0.9725758139465345
$ws.Range("P5").Value = 0.9725758139465346
$ws.Range("Q5").Value = 26.03031732377667
$ws.Range("R5").Value = 234.27285591399
$ws.Range("S5").Value = 0.07200863595348427
$ws.Range("T5").Value = 0.07200863595348428

# Row 6: FAPs -> MuSCs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Trf"
$ws.Range("C6").Value = "Tfr2"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.234235
$ws.Range("H6").Value = 27.702705
$ws.Range("I6").Value = 0.07403909795092109
$ws.Range("J6").Value = 0.07403909795092109
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.04134966666666667
$ws.Range("N6").Value = 0.124049
$ws.Range("O6").Value = 0.01426648349910611
$ws.Range("P6").Value = 0.01426648349910611
$ws.Range("Q6").Value = 0.3818325391716667
$ws.Range("R6").Value = 3.436492852545
$ws.Range("S6").Value = 0.001056277569205517
$ws.Range("T6").Value = 0.001056277569205517

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Trf"
$ws.Range("C7").Value = "Tfr2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.234235
$ws.Range("H7").Value = 27.702705
$ws.Range("I7").Value = 0.07403909795092109
$ws.Range("J7").Value = 0.07403909795092109
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.038136
$ws.Range("N7").Value = 0.114408
$ws.Range("O7").Value = 0.01315770255435942
$ws.Range("P7").Value = 0.01315770255435942
$ws.Range("Q7").Value = 0.3521567859599999
$ws.Range("R7").Value = 3.16941107364
$ws.Range("S7").Value = 0.0009741844282313016
$ws.Range("T7").Value = 0.0009741844282313019

# Row 8: MuSCs -> FAPs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Trf"
$ws.Range("C8").Value = "Tfr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.276267
$ws.Range("H8").Value = 3.828801
$ws.Range("I8").Value = 0.01023297083348303
$ws.Range("J8").Value = 0.01023297083348304
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.818892666666667
$ws.Range("N8").Value = 8.456678
$ws.Range("O8").Value = 0.9725758139465345
$ws.Range("P8").Value = 0.9725758139465346
$ws.Range("Q8").Value = 3.597659687008667
$ws.Range("R8").Value = 32.378937183078
$ws.Range("S8").Value = 0.009952339937465908
$ws.Range("T8").Value = 0.009952339937465911

# Row 9: MuSCs -> MuSCs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Trf"
$ws.Range("C9").Value = "Tfr2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.276267
$ws.Range("H9").Value = 3.828801
$ws.Range("I9").Value = 0.01023297083348303
$ws.Range("J9").Value = 0.01023297083348304
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.04134966666666667
$ws.Range("N9").Value = 0.124049
$ws.Range("O9").Value = 0.01426648349910611
$ws.Range("P9").Value = 0.01426648349910611
$ws.Range("Q9").Value = 0.05277321502766667
$ws.Range("R9").Value = 0.4749589352490001
$ws.Range("S9").Value = 0.0001459885095427198
$ws.Range("T9").Value = 0.0001459885095427198

# Row 10: MuSCs -> Resolving-Mac
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Trf"
$ws.Range("C10").Value = "Tfr2"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.276267
$ws.Range("H10").Value = 3.828801
$ws.Range("I10").Value = 0.01023297083348303
$ws.Range("J10").Value = 0.01023297083348304
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.038136
$ws.Range("N10").Value = 0.114408
$ws.Range("O10").Value = 0.01315770255435942
$ws.Range("P10").Value = 0.01315770255435942
$ws.Range("Q10").Value = 0.048671718312
$ws.Range("R10").Value = 0.438045464808
$ws.Range("S10").Value = 0.0001346423864744051
$ws.Range("T10").Value = 0.0001346423864744052

# Row 11: Resolving-Mac -> FAPs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Trf"
$ws.Range("C11").Value = "Tfr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 112.6100793333333
$ws.Range("H11").Value = 337.830238
$ws.Range("I11").Value = 0.902895442234431
$ws.Range("J11").Value = 0.9028954422344311
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.818892666666667
$ws.Range("N11").Value = 8.456678
$ws.Range("O11").Value = 0.9725758139465345
$ws.Range("P11").Value = 0.9725758139465346
$ws.Range("Q11").Value = 317.4357268254849
$ws.Range("R11").Value = 2856.921541429364
$ws.Range("S11").Value = 0.8781342696397679
$ws.Range("T11").Value = 0.8781342696397682

# Row 12: Resolving-Mac -> MuSCs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Trf"
$ws.Range("C12").Value = "Tfr2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 112.6100793333333
$ws.Range("H12").Value = 337.830238
$ws.Range("I12").Value = 0.902895442234431
$ws.Range("J12").Value = 0.9028954422344311
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.04134966666666667
$ws.Range("N12").Value = 0.124049
$ws.Range("O12").Value = 0.01426648349910611
$ws.Range("P12").Value = 0.01426648349910611
$ws.Range("Q12").Value = 4.656389243740223
$ws.Range("R12").Value = 41.907503193662
$ws.Range("S12").Value = 0.01288114292805562
$ws.Range("T12").Value = 0.01288114292805563

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Trf"
$ws.Range("C13").Value = "Tfr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 112.6100793333333
$ws.Range("H13").Value = 337.830238
$ws.Range("I13").Value = 0.902895442234431
$ws.Range("J13").Value = 0.9028954422344311
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.038136
$ws.Range("N13").Value = 0.114408
$ws.Range("O13").Value = 0.01315770255435942
$ws.Range("P13").Value = 0.01315770255435942
$ws.Range("Q13").Value = 4.294497985455999
$ws.Range("R13").Value = 38.650481869104
$ws.Range("S13").Value = 0.01188002966660745
$ws.Range("T13").Value = 0.01188002966660745
